$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.064.40"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.247.04"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.36"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.50"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.522"
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.69"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").Value = "2.597.45"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.43"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "2.256.98"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.777"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").Value = "41.914.34"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.13"
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("D20").Value = "0.0₃0900"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.90"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.08"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.17"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.33"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.01"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.46"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.03"
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.15"
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.51"
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0718"
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.05"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").Value = "1.943.81"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0280"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.39"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.17"
$ws.Range("E45").Value = "  -9.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("E47").Value = "  -3.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.51"
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("D49").Value = "2.468.74"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.08"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.86"
$ws.Range("E51").Value = "  -1.17%  "

Write-Host "Applied cryptos update"
